$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header row: "_old" suffix -> "_FV2404", "_new" suffix -> "_FV2410"
#    Columns A-J hold the "_old" headers, column K is "diff" (unchanged),
#    columns L-U hold the "_new" headers.
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = $cell.Value2.ToString().Replace("_old", "_FV2404")
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = $cell.Value2.ToString().Replace("_new", "_FV2410")
}

# 2) Turn the populated range into an Excel Table ("Table1") covering A1:U76
$rng = $ws.Range("A1:U76")
$lo = $ws.ListObjects.Add(1, $rng, 0, 1)
$lo.Name = "Table1"

# 3) Freeze the header row (split after row 1, frozen)
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Output "Headers renamed, Table1 created (A1:U76), header row frozen."
